# Apply scraped-schedule refresh to workbook (run at 08:00:50)
$wb = $excel.ActiveWorkbook

# ---- Sheet: LP1912 ----
$ws = $wb.Worksheets.Item("LP1912")
$ws.Cells.Item(2,1).Value = "Última actualización: 08:00:50"
$ws.Cells.Item(3,1).Value = "Total filas: 90"
$ws.Range("A59:E85").ClearContents()
$data = @(
    @(59,"08:00:50","08:00","11_ETCHEVERRY",0,"LP1912"),
    @(60,"08:00:50","08:00","215A_EL PATO",0,"LP1912"),
    @(61,"07:48:14","08:01","16_SANTA ANA",13,"LP1912"),
    @(62,"08:00:50","08:03","17X38_ROMERO",3,"LP1912"),
    @(63,"06:46:06","08:03","23_HERNANDEZ",77,"LP1912"),
    @(64,"06:58:01","08:04","23_HERNANDEZ",66,"LP1912"),
    @(65,"08:00:50","08:14","10_OLMOS",14,"LP1912"),
    @(66,"08:00:50","08:19","15_ABASTO",19,"LP1912"),
    @(67,"08:00:50","08:21","16_SANTA ANA",21,"LP1912"),
    @(68,"08:00:50","08:29","14_ABASTO",29,"LP1912"),
    @(69,"07:48:14","08:30","14_ABASTO",42,"LP1912"),
    @(70,"08:00:50","08:33","23_HERNANDEZ",33,"LP1912"),
    @(71,"08:00:50","08:33","215C_EL PATO",33,"LP1912"),
    @(72,"07:48:14","08:34","215C_EL PATO",46,"LP1912"),
    @(73,"08:00:50","08:41","16_SANTA ANA",41,"LP1912"),
    @(74,"08:00:50","08:44","10_OLMOS",44,"LP1912"),
    @(75,"08:00:50","08:47","215A_EL PATO",47,"LP1912"),
    @(76,"07:48:14","08:48","215A_EL PATO",60,"LP1912"),
    @(77,"07:21:42","08:51","16_P MOR-SANTA ANA",90,"LP1912"),
    @(78,"07:48:14","08:52","16_P MOR-SANTA ANA",64,"LP1912"),
    @(79,"08:00:50","08:55","16_P MOR-SANTA ANA",55,"LP1912"),
    @(80,"08:00:50","08:59","215B_EL PATO",59,"LP1912"),
    @(81,"07:48:14","09:00","215B_EL PATO",72,"LP1912"),
    @(82,"08:00:50","09:03","17X38_ROMERO",63,"LP1912"),
    @(83,"08:00:50","09:03","23_HERNANDEZ",63,"LP1912"),
    @(84,"08:00:50","09:14","27_EL RETIRO",74,"LP1912"),
    @(85,"08:00:50","09:14","11_ETCHEVERRY",74,"LP1912"),
    @(86,"07:48:14","09:15","11_ETCHEVERRY",87,"LP1912"),
    @(87,"07:48:14","09:16","27_EL RETIRO",88,"LP1912"),
    @(88,"08:00:50","09:18","215_EL PELIGRO",78,"LP1912"),
    @(89,"07:48:14","09:19","215_EL PELIGRO",91,"LP1912"),
    @(90,"07:48:14","09:29","10_OLMOS",101,"LP1912"),
    @(91,"08:00:50","09:34","15_ABASTO",94,"LP1912"),
    @(92,"08:00:50","09:44","14_ABASTO",104,"LP1912"),
    @(93,"07:48:14","09:45","14_ABASTO",117,"LP1912"),
    @(94,"08:00:50","09:51","16_P MOR-SANTA ANA",111,"LP1912"),
    @(95,"08:00:50","09:56","10_OLMOS",116,"LP1912")
)
foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r,1).Value = $row[1]
    $ws.Cells.Item($r,2).Value = $row[2]
    $ws.Cells.Item($r,3).Value = $row[3]
    $ws.Cells.Item($r,4).Value = $row[4]
    $ws.Cells.Item($r,5).Value = $row[5]
}

# ---- Sheet: LP1912-215 ----
$ws = $wb.Worksheets.Item("LP1912-215")
$ws.Cells.Item(2,1).Value = "Última actualización: 08:00:50"
$ws.Cells.Item(3,1).Value = "Total filas: 20"
$ws.Range("A17:E23").ClearContents()
$data = @(
    @(17,"08:00:50","08:00","215A_EL PATO",0,"LP1912"),
    @(18,"08:00:50","08:33","215C_EL PATO",33,"LP1912"),
    @(19,"07:48:14","08:34","215C_EL PATO",46,"LP1912"),
    @(20,"08:00:50","08:47","215A_EL PATO",47,"LP1912"),
    @(21,"07:48:14","08:48","215A_EL PATO",60,"LP1912"),
    @(22,"08:00:50","08:59","215B_EL PATO",59,"LP1912"),
    @(23,"07:48:14","09:00","215B_EL PATO",72,"LP1912"),
    @(24,"08:00:50","09:18","215_EL PELIGRO",78,"LP1912"),
    @(25,"07:48:14","09:19","215_EL PELIGRO",91,"LP1912")
)
foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r,1).Value = $row[1]
    $ws.Cells.Item($r,2).Value = $row[2]
    $ws.Cells.Item($r,3).Value = $row[3]
    $ws.Cells.Item($r,4).Value = $row[4]
    $ws.Cells.Item($r,5).Value = $row[5]
}

# ---- Sheet: 6203-6173 ----
$ws = $wb.Worksheets.Item("6203-6173")
$ws.Cells.Item(2,1).Value = "Última actualización: 08:00:50"
$ws.Cells.Item(3,1).Value = "Total filas: 15"
$data = @(
    @(13,"08:00:50","08:09","215A_LA PLATA",9,"L6173"),
    @(16,"08:00:50","08:26","215C_LA PLATA",26,"L6203"),
    @(19,"08:00:50","09:20","215A_LA PLATA",80,"L6173"),
    @(20,"08:00:50","09:55","215C_LA PLATA",115,"L6203")
)
foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r,1).Value = $row[1]
    $ws.Cells.Item($r,2).Value = $row[2]
    $ws.Cells.Item($r,3).Value = $row[3]
    $ws.Cells.Item($r,4).Value = $row[4]
    $ws.Cells.Item($r,5).Value = $row[5]
}
